# Natmi following Dr Hou advice
# Efna5-Epha3 LR-pairs results recomputed after adding the "ECs" cluster.
# This rewrites the two existing result rows and appends four new rows so
# that every Sending-cluster x Target-cluster combination of
# {FAPs, sCs} x {ECs, FAPs, sCs} is present.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = @("FAPs", "Efna5", "Epha3", "ECs", 3, 1, 1.666083666666667, 4.998251, 0.6125276070882968, 0.6125276070882968, 1, 0.3333333333333333, 0.06698166666666668, 0.200945, 0.003012576978541733, 0.003012576978541732, 0.1115970607994445, 1.004373547195, 0.001845286567835459, 0.001845286567835459)
$row3 = @("FAPs", "Efna5", "Epha3", "FAPs", 3, 1, 1.666083666666667, 4.998251, 0.6125276070882968, 0.6125276070882968, 3, 1, 22.03620333333333, 66.10861, 0.9911034191912899, 0.9911034191912899, 36.71415844901222, 330.42742604111, 0.6070782057342698, 0.6070782057342698)
$row4 = @("FAPs", "Efna5", "Epha3", "sCs", 3, 1, 1.666083666666667, 4.998251, 0.6125276070882968, 0.6125276070882968, 2, 0.6666666666666666, 0.130825, 0.392475, 0.005884003830168287, 0.005884003830168287, 0.2179653956916666, 1.961688561225, 0.003604114786191354, 0.003604114786191354)
$row5 = @("sCs", "Efna5", "Epha3", "ECs", 3, 1, 1.053930333333333, 3.161791, 0.3874723929117032, 0.3874723929117031, 1, 0.3333333333333333, 0.06698166666666668, 0.200945, 0.003012576978541733, 0.003012576978541732, 0.07059401027722223, 0.6353460924950001, 0.001167290410706274, 0.001167290410706274)
$row6 = @("sCs", "Efna5", "Epha3", "FAPs", 3, 1, 1.053930333333333, 3.161791, 0.3874723929117032, 0.3874723929117031, 3, 1, 22.03620333333333, 66.10861, 0.9911034191912899, 0.9911034191912899, 23.22462312450111, 209.02160812051, 0.38402521345702, 0.3840252134570198)
$row7 = @("sCs", "Efna5", "Epha3", "sCs", 3, 1, 1.053930333333333, 3.161791, 0.3874723929117032, 0.3874723929117031, 2, 0.6666666666666666, 0.130825, 0.392475, 0.005884003830168287, 0.005884003830168287, 0.1378804358583333, 1.240923922725, 0.002279889043976933, 0.002279889043976933)

$data = @($row2, $row3, $row4, $row5, $row6, $row7)

# Rows are written starting at row 2 (row 1 holds the column headers),
# column A is column 1, so $thisRow[$c] -> column ($c + 1).
for ($r = 0; $r -lt $data.Length; $r++) {
    $thisRow = $data[$r]
    for ($c = 0; $c -lt $thisRow.Length; $c++) {
        $ws.Cells.Item(2 + $r, $c + 1).Value2 = $thisRow[$c]
    }
}
